$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 88
$ws.Range("B7").Value = "ew90"
$ws.Range("B8").Value = "pranav"

$ws.Range("B9").Select()
